# The workbook contains four repeated 8-row blocks (rows 2-9, 10-17, 18-25,
# 26-33) reporting subunit-count stats for four different GO groups:
#   rows 2-9:   "intrinsic component of membrane"
#   rows 10-17: "periplasmic space"
#   rows 18-25: "plasma membrane"
#   rows 26-33: "cytosol"
#
# This edit removes the "periplasmic space" block entirely (rows 10-17).
# Deleting those rows shifts the "plasma membrane" and "cytosol" blocks up
# by 8 rows (to 10-17 and 18-25 respectively) and drops the now-unused
# "periplasmic space" shared string, which Excel automatically does when a
# whole-row delete leaves no remaining references to that string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:17").Delete()
